$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from an existing data row (row 2) onto the new row 5
# so the new cells get the same cell style ("s=1") as the rest of the table.
$ws.Range("A2:C2").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add the new row of data
$ws.Range("A5").Value = "canopus"
$ws.Range("B5").Value = "Darllan"
$ws.Range("C5").Value = "darllan.dias@estasa.com.br"
